$wb = $excel.ActiveWorkbook

# --- Sheet references ---
$wsCuenta = $wb.Worksheets.Item("DatosCuenta")
$wsHogar  = $wb.Worksheets.Item("DatosHogar")
$wsMotor  = $wb.Worksheets.Item("DatosMotor")
$wsAP     = $wb.Worksheets.Item("DatosAP")

# --- DatosCuenta (sheet1) updates ---
$wsCuenta.Range("A2").Value = "SmokeName"
$wsCuenta.Range("B2").Value = "SmokeLastName"
$wsCuenta.Range("C2").Value = 23100200
$wsCuenta.Range("D2").Value = 123

# --- DatosHogar (sheet2) updates ---
$wsHogar.Range("A2").Value = 640

# --- DatosMotor (sheet3) updates ---
$wsMotor.Range("A2").Value = "SQA040"
$wsMotor.Range("B2").Value = "ABC12SSQA040"
$wsMotor.Range("C2").Value = "ZAZ123SSQA040"

# --- DatosAP (sheet4) updates ---
$wsAP.Range("A2").Value = 21200140

# --- Selections on each sheet ---
$wsCuenta.Range("G11").Select()
$wsMotor.Range("D10").Select()
$wsAP.Range("F11").Select()

# --- Activate DatosAP last so it becomes the active/selected tab ---
$wsAP.Activate()
$wsAP.Range("F11").Select()
